# Apply cryptos list update (cell text/value changes only; no formatting changes).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$sub3 = [string][char]0x2083   # subscript 3, used in "0.0<sub3>0702"
$sub6 = [string][char]0x2086   # subscript 6, used in "0.0<sub6>0104"

$ws.Cells.Item(2, 4).Value = "27.203.25"
$ws.Cells.Item(2, 5).Value = "  +0.80%  "
$ws.Cells.Item(3, 4).Value = "1.570.33"
$ws.Cells.Item(3, 5).Value = "  +0.66%  "
$ws.Cells.Item(4, 5).Value = "  +0.47%  "
$ws.Cells.Item(5, 4).NumberFormat = "@"
$ws.Cells.Item(5, 4).Value = "211.66"
$ws.Cells.Item(5, 4).Style = "Normal"
$ws.Cells.Item(5, 5).Value = "  +2.05%  "
$ws.Cells.Item(6, 5).Value = "  +0.74%  "
$ws.Cells.Item(7, 5).Value = "  +0.48%  "
$ws.Cells.Item(8, 4).NumberFormat = "@"
$ws.Cells.Item(8, 4).Value = "22.10"
$ws.Cells.Item(8, 4).Style = "Normal"
$ws.Cells.Item(8, 5).Value = "  -0.10%  "
$ws.Cells.Item(9, 5).Value = "  +0.48%  "
$ws.Cells.Item(10, 4).NumberFormat = "@"
$ws.Cells.Item(10, 4).Value = "0.0600"
$ws.Cells.Item(10, 4).Style = "Normal"
$ws.Cells.Item(10, 5).Value = "  +0.65%  "
$ws.Cells.Item(11, 4).NumberFormat = "@"
$ws.Cells.Item(11, 4).Value = "0.0867"
$ws.Cells.Item(11, 4).Style = "Normal"
$ws.Cells.Item(11, 5).Value = "  +0.94%  "
$ws.Cells.Item(12, 4).Value = "1.793.73"
$ws.Cells.Item(12, 5).Value = "  +0.63%  "
$ws.Cells.Item(13, 4).Value = "1.575.64"
$ws.Cells.Item(13, 5).Value = "  +1.25%  "
$ws.Cells.Item(14, 5).Value = "  +0.81%  "
$ws.Cells.Item(15, 5).Value = "  +0.03%  "
$ws.Cells.Item(16, 4).Value = "27.195.29"
$ws.Cells.Item(16, 5).Value = "  +0.74%  "
$ws.Cells.Item(17, 4).NumberFormat = "@"
$ws.Cells.Item(17, 4).Value = "62.31"
$ws.Cells.Item(17, 4).Style = "Normal"
$ws.Cells.Item(17, 5).Value = "  +0.32%  "
$ws.Cells.Item(18, 4).Value = "0.0" + $sub3 + "0702"
$ws.Cells.Item(18, 5).Value = "  -0.34%  "
$ws.Cells.Item(19, 2).Value = "Chainlink"
$ws.Cells.Item(19, 3).Value = "https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link"
$ws.Cells.Item(19, 4).NumberFormat = "@"
$ws.Cells.Item(19, 4).Value = "7.43"
$ws.Cells.Item(19, 4).Style = "Normal"
$ws.Cells.Item(19, 5).Value = "  +0.81%  "
$ws.Cells.Item(20, 2).Value = "BitcoinCash"
$ws.Cells.Item(20, 3).Value = "https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch"
$ws.Cells.Item(20, 4).NumberFormat = "@"
$ws.Cells.Item(20, 4).Value = "216.29"
$ws.Cells.Item(20, 4).Style = "Normal"
$ws.Cells.Item(20, 5).Value = "  -0.43%  "
$ws.Cells.Item(21, 5).Value = "  +0.47%  "
$ws.Cells.Item(22, 5).Value = "  +1.10%  "
$ws.Cells.Item(23, 4).NumberFormat = "@"
$ws.Cells.Item(23, 4).Value = "9.24"
$ws.Cells.Item(23, 4).Style = "Normal"
$ws.Cells.Item(23, 5).Value = "  -0.19%  "
$ws.Cells.Item(24, 5).Value = "  +1.05%  "
$ws.Cells.Item(25, 4).NumberFormat = "@"
$ws.Cells.Item(25, 4).Value = "153.99"
$ws.Cells.Item(25, 4).Style = "Normal"
$ws.Cells.Item(25, 5).Value = "  +0.48%  "
$ws.Cells.Item(26, 4).NumberFormat = "@"
$ws.Cells.Item(26, 4).Value = "6.70"
$ws.Cells.Item(26, 4).Style = "Normal"
$ws.Cells.Item(26, 5).Value = "  +1.22%  "
$ws.Cells.Item(27, 4).NumberFormat = "@"
$ws.Cells.Item(27, 4).Value = "15.09"
$ws.Cells.Item(27, 4).Style = "Normal"
$ws.Cells.Item(27, 5).Value = "  +0.18%  "
$ws.Cells.Item(28, 5).Value = "  +2.40%  "
$ws.Cells.Item(29, 5).Value = "  +0.42%  "
$ws.Cells.Item(30, 5).Value = "  +2.15%  "
$ws.Cells.Item(31, 5).Value = "  +1.24%  "
$ws.Cells.Item(32, 4).NumberFormat = "@"
$ws.Cells.Item(32, 4).Value = "3.24"
$ws.Cells.Item(32, 4).Style = "Normal"
$ws.Cells.Item(32, 5).Value = "  +0.21%  "
$ws.Cells.Item(33, 4).Value = "1.450.90"
$ws.Cells.Item(33, 5).Value = "  +2.02%  "
$ws.Cells.Item(34, 5).Value = "  +2.02%  "
$ws.Cells.Item(35, 5).Value = "  +4.91%  "
$ws.Cells.Item(36, 5).Value = "  +0.34%  "
$ws.Cells.Item(37, 4).NumberFormat = "@"
$ws.Cells.Item(37, 4).Value = "2.35"
$ws.Cells.Item(37, 4).Style = "Normal"
$ws.Cells.Item(37, 5).Value = "  +1.24%  "
$ws.Cells.Item(38, 5).Value = "  +0.92%  "
$ws.Cells.Item(39, 4).NumberFormat = "@"
$ws.Cells.Item(39, 4).Value = "0.537"
$ws.Cells.Item(39, 4).Style = "Normal"
$ws.Cells.Item(39, 5).Value = "  +0.66%  "
$ws.Cells.Item(40, 5).Value = "  +2.33%  "
$ws.Cells.Item(41, 4).NumberFormat = "@"
$ws.Cells.Item(41, 4).Value = "0.809"
$ws.Cells.Item(41, 4).Style = "Normal"
$ws.Cells.Item(41, 5).Value = "  +0.02%  "
$ws.Cells.Item(42, 5).Value = "  +0.43%  "
$ws.Cells.Item(43, 5).Value = "  +0.77%  "
$ws.Cells.Item(44, 5).Value = "  +0.23%  "
$ws.Cells.Item(45, 4).NumberFormat = "@"
$ws.Cells.Item(45, 4).Value = "64.66"
$ws.Cells.Item(45, 4).Style = "Normal"
$ws.Cells.Item(45, 5).Value = "  -0.32%  "
$ws.Cells.Item(46, 4).NumberFormat = "@"
$ws.Cells.Item(46, 4).Value = "1.73"
$ws.Cells.Item(46, 4).Style = "Normal"
$ws.Cells.Item(46, 5).Value = "  -0.61%  "
$ws.Cells.Item(47, 4).Value = "1.706.26"
$ws.Cells.Item(47, 5).Value = "  +0.63%  "
$ws.Cells.Item(48, 4).NumberFormat = "@"
$ws.Cells.Item(48, 4).Value = "85.99"
$ws.Cells.Item(48, 4).Style = "Normal"
$ws.Cells.Item(48, 5).Value = "  -1.69%  "
$ws.Cells.Item(49, 4).Value = "0.0" + $sub6 + "0104"
$ws.Cells.Item(49, 5).Value = "  +3.87%  "
$ws.Cells.Item(50, 4).NumberFormat = "@"
$ws.Cells.Item(50, 4).Value = "0.0520"
$ws.Cells.Item(50, 4).Style = "Normal"
$ws.Cells.Item(51, 4).NumberFormat = "@"
$ws.Cells.Item(51, 4).Value = "0.0961"
$ws.Cells.Item(51, 4).Style = "Normal"
$ws.Cells.Item(51, 5).Value = "  +0.58%  "
